# Add new columns I (I0) and J (IF) to the worksheet, mirroring the
# formatting of the existing header row and copying values from column H
# (with one exception on row 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# Match the look of the other header cells (bold / centered / bordered style)
# by copying the formatting from H1, then set the new header text.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)
$ws.Range("J1").Value = "IF"

$excel.CutCopyMode = $false

# --- Data rows (2-25) ---
# Column I is 1 for every row except row 22 (which is 7).
# Column J mirrors column H for every row except row 22 (which is H22 + 6 = 8).
for ($r = 2; $r -le 25; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()

    if ($r -eq 22) {
        $iVal = 7
        $jVal = 8
    } else {
        $iVal = 1
        $jVal = $hVal
    }

    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
